$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '27.681.35'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.40%  '
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.637.85'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -0.68%  '
$c.Style = "Normal"
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '212.50'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.523'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -2.16%  '
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '23.03'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -2.55%  '
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +0.53%  '
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.Style = "Normal"
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.869.74'
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.656.44'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +0.51%  '
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -4.32%  '
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '64.63'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.06%  '
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '27.654.88'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '230.18'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -0.65%  '
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.70'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +1.87%  '
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.0₃0724'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.30'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '10.23'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +4.29%  '
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +1.83%  '
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '150.65'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +1.14%  '
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '6.95'
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -1.01%  '
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -1.58%  '
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '15.59'
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +0.38%  '
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.0486'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -0.32%  '
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.457.10'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +1.87%  '
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '3.12'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -2.53%  '
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -2.41%  '
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.54%  '
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -0.86%  '
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.900'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +10.01%  '
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '69.11'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +5.84%  '
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -1.32%  '
$c.Style = "Normal"
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '5.57'
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +0.54%  '
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -0.97%  '
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.779.62'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.78%  '
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +2.93%  '
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '86.80'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.Style = "Normal"
$c = $ws.Range('B50')
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.Style = "Normal"
$c = $ws.Range('C50')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0993'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"
$c = $ws.Range('B51')
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c.Style = "Normal"
$c = $ws.Range('C51')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.80'
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.Style = "Normal"
